$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.725.24'
$ws.Range("E2").Value = '  +1.17%  '

$ws.Range("D3").Value = '3.350.78'
$ws.Range("E3").Value = '  +1.09%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '587.24'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +5.61%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '187.18'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.66%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.603'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +3.67%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("E9").Value = '  +1.45%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.587'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.65%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '47.32'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.75%  '

$ws.Range("E12").Value = '  +1.65%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '650.66'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +7.79%  '

$ws.Range("D14").Value = '3.884.72'
$ws.Range("E14").Value = '  +1.08%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '8.56'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.01%  '

$ws.Range("D16").Value = '66.687.65'
$ws.Range("E16").Value = '  +1.09%  '

$ws.Range("E17").Value = '  +0.65%  '

$ws.Range("D18").Value = '3.346.86'
$ws.Range("E18").Value = '  +1.08%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '17.94'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.16%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '11.17'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.07%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.905'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.44%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '17.78'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -4.99%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.11'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '100.87'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.24%  '

$ws.Range("E25").Value = '  +1.59%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.82'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +2.71%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.68'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.71%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '32.10'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +6.03%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '8.65'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.47%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.94'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +2.91%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '608.07'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +5.08%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.86'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.28%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '11.20'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.49%  '

$ws.Range("D34").Value = '3.887.49'
$ws.Range("E34").Value = '  +4.82%  '

$ws.Range("E35").Value = '  +1.52%  '

$ws.Range("E36").Value = '  -0.02%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '55.78'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.19%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.77'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +4.08%  '

$ws.Range("E39").Value = '  +1.78%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '33.65'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.38%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.22'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.25%  '

$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").Value = '0.0₃0702'
$ws.Range("E42").Value = '  -0.19%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.343'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.31%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.37'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.10%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0419'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.38%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.131'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.32%  '

$ws.Range("E47").Value = '  +0.09%  '

$ws.Range("E48").Value = '  +0.31%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.36'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +9.05%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.84'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -18.00%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '130.71'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +5.97%  '
